{"js": "// Replace each three-digit x one-digit multiplication expression with its\n// updated value. Every \"old\" expression in this document is unique, so a\n// scoped search-and-replace by exact text is unambiguous.\nconst replacements = [\n  [\"646\u00d79=5814\", \"914\u00d72=1828\"],\n  [\"810\u00d73=2430\", \"138\u00d78=1104\"],\n  [\"219\u00d73=657\", \"994\u00d75=4970\"],\n  [\"171\u00d73=513\", \"344\u00d77=2408\"],\n  [\"314\u00d73=942\", \"253\u00d76=1518\"],\n  [\"379\u00d79=3411\", \"298\u00d78=2384\"],\n  [\"922\u00d77=6454\", \"359\u00d76=2154\"],\n  [\"363\u00d77=2541\", \"510\u00d72=1020\"],\n  [\"641\u00d76=3846\", \"132\u00d78=1056\"],\n  [\"930\u00d74=3720\", \"965\u00d75=4825\"],\n  [\"958\u00d74=3832\", \"264\u00d76=1584\"],\n  [\"852\u00d76=5112\", \"815\u00d75=4075\"],\n  [\"949\u00d74=3796\", \"401\u00d75=2005\"],\n  [\"422\u00d79=3798\", \"518\u00d77=3626\"],\n  [\"280\u00d74=1120\", \"746\u00d78=5968\"],\n  [\"254\u00d75=1270\", \"656\u00d79=5904\"],\n  [\"364\u00d74=1456\", \"779\u00d76=4674\"],\n  [\"987\u00d77=6909\", \"947\u00d78=7576\"],\n  [\"267\u00d79=2403\", \"576\u00d79=5184\"],\n  [\"495\u00d76=2970\", \"525\u00d79=4725\"],\n  [\"238\u00d77=1666\", \"895\u00d75=4475\"],\n  [\"364\u00d72=728\", \"348\u00d74=1392\"],\n  [\"945\u00d77=6615\", \"380\u00d78=3040\"],\n  [\"460\u00d75=2300\", \"570\u00d75=2850\"],\n  [\"756\u00d77=5292\", \"913\u00d79=8217\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit x one-digit multiplication expression with its\n# updated value. Every \"old\" expression in this document is unique, so a\n# document-wide Find/Replace by exact text is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"646\u00d79=5814\", \"914\u00d72=1828\"),\n    @(\"810\u00d73=2430\", \"138\u00d78=1104\"),\n    @(\"219\u00d73=657\", \"994\u00d75=4970\"),\n    @(\"171\u00d73=513\", \"344\u00d77=2408\"),\n    @(\"314\u00d73=942\", \"253\u00d76=1518\"),\n    @(\"379\u00d79=3411\", \"298\u00d78=2384\"),\n    @(\"922\u00d77=6454\", \"359\u00d76=2154\"),\n    @(\"363\u00d77=2541\", \"510\u00d72=1020\"),\n    @(\"641\u00d76=3846\", \"132\u00d78=1056\"),\n    @(\"930\u00d74=3720\", \"965\u00d75=4825\"),\n    @(\"958\u00d74=3832\", \"264\u00d76=1584\"),\n    @(\"852\u00d76=5112\", \"815\u00d75=4075\"),\n    @(\"949\u00d74=3796\", \"401\u00d75=2005\"),\n    @(\"422\u00d79=3798\", \"518\u00d77=3626\"),\n    @(\"280\u00d74=1120\", \"746\u00d78=5968\"),\n    @(\"254\u00d75=1270\", \"656\u00d79=5904\"),\n    @(\"364\u00d74=1456\", \"779\u00d76=4674\"),\n    @(\"987\u00d77=6909\", \"947\u00d78=7576\"),\n    @(\"267\u00d79=2403\", \"576\u00d79=5184\"),\n    @(\"495\u00d76=2970\", \"525\u00d79=4725\"),\n    @(\"238\u00d77=1666\", \"895\u00d75=4475\"),\n    @(\"364\u00d72=728\", \"348\u00d74=1392\"),\n    @(\"945\u00d77=6615\", \"380\u00d78=3040\"),\n    @(\"460\u00d75=2300\", \"570\u00d75=2850\"),\n    @(\"756\u00d77=5292\", \"913\u00d79=8217\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
